$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.243.32'
$ws.Range("E2").Value = '  -1.99%  '
$ws.Range("D3").Value = '2.581.61'
$ws.Range("E3").Value = '  -2.15%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = "'563.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.16%  '
$ws.Range("D6").Value = "'142.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.63%  '
$ws.Range("E7").Value = '  +0.30%  '
$ws.Range("E8").Value = '  -2.24%  '
$ws.Range("D9").Value = '2.588.26'
$ws.Range("E9").Value = '  -2.68%  '
$ws.Range("D10").Value = "'6.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.80%  '
$ws.Range("E11").Value = '  -0.64%  '
$ws.Range("E12").Value = '  +11.53%  '
$ws.Range("E13").Value = '  +3.05%  '
$ws.Range("D14").Value = '3.034.44'
$ws.Range("E14").Value = '  -2.68%  '
$ws.Range("D15").Value = "'23.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +7.27%  '
$ws.Range("D16").Value = '59.200.48'
$ws.Range("E16").Value = '  -2.00%  '
$ws.Range("E17").Value = '  +0.31%  '
$ws.Range("D18").Value = '2.585.60'
$ws.Range("E18").Value = '  -2.81%  '
$ws.Range("D19").Value = "'4.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.62%  '
$ws.Range("D20").Value = "'337.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.21%  '
$ws.Range("D21").Value = "'10.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.64%  '
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("E24").Value = '  -4.06%  '
$ws.Range("E25").Value = '  +5.05%  '
$ws.Range("E26").Value = '  +0.30%  '
$ws.Range("E27").Value = '  -3.13%  '
$ws.Range("E28").Value = '  -0.36%  '
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("E30").Value = '  +0.07%  '
$ws.Range("E31").Value = '  -2.56%  '
$ws.Range("D32").Value = "'160.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.51%  '
$ws.Range("E33").Value = '  -0.12%  '
$ws.Range("E34").Value = '  -1.23%  '
$ws.Range("D35").Value = "'4.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.35%  '
$ws.Range("E36").Value = '  -0.75%  '
$ws.Range("E37").Value = '  -3.46%  '
$ws.Range("E38").Value = '  -3.63%  '
$ws.Range("E39").Value = '  -0.43%  '
$ws.Range("E40").Value = '  -2.32%  '
$ws.Range("D41").Value = "'294.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.49%  '
$ws.Range("D42").Value = "'3.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("E43").Value = '  +0.56%  '
$ws.Range("D44").Value = "'132.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.72%  '
$ws.Range("E45").Value = '  -1.02%  '
$ws.Range("E46").Value = '  -1.86%  '
$ws.Range("E47").Value = '  -0.10%  '
$ws.Range("E48").Value = '  -2.41%  '
$ws.Range("E49").Value = '  -2.36%  '
$ws.Range("E50").Value = '  -0.93%  '
$ws.Range("D51").Value = "'18.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.83%  '
